function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws 'D2' '26.593.19'
Set-TextValue $ws 'D3' '1.692.09'
Set-TextValue $ws 'E3' '  -5.64%  '
Set-TextValue $ws 'E4' '  +0.23%  '
Set-TextValue $ws 'D5' '219.78'
Set-TextValue $ws 'E5' '  -5.08%  '
Set-TextValue $ws 'D6' '0.5110'
Set-TextValue $ws 'E6' '  -13.09%  '
Set-TextValue $ws 'E7' '  +0.16%  '
Set-TextValue $ws 'D8' '0.2659'
Set-TextValue $ws 'E8' '  -3.93%  '
Set-TextValue $ws 'D9' '22.06'
Set-TextValue $ws 'E9' '  -4.75%  '
Set-TextValue $ws 'D10' '0.06323'
Set-TextValue $ws 'E10' '  -6.19%  '
Set-TextValue $ws 'D11' '0.07364'
Set-TextValue $ws 'E11' '  -2.21%  '
Set-TextValue $ws 'D12' '1.697.88'
Set-TextValue $ws 'E12' '  -5.34%  '
Set-TextValue $ws 'E13' '  -5.80%  '
Set-TextValue $ws 'D14' '0.5784'
Set-TextValue $ws 'E14' '  -5.70%  '
Set-TextValue $ws 'D15' '1.922.73'
Set-TextValue $ws 'E15' '  -5.58%  '
Set-TextValue $ws 'D16' '0.000008512'
Set-TextValue $ws 'E16' '  -5.53%  '
Set-TextValue $ws 'D17' '65.36'
Set-TextValue $ws 'E17' '  -13.20%  '
Set-TextValue $ws 'D18' '26.618.02'
Set-TextValue $ws 'E18' '  -6.93%  '
Set-TextValue $ws 'D19' '4.982'
Set-TextValue $ws 'E19' '  -8.20%  '
Set-TextValue $ws 'D20' '1.005'
Set-TextValue $ws 'E20' '  +0.12%  '
Set-TextValue $ws 'D21' '10.95'
Set-TextValue $ws 'E21' '  -4.52%  '
Set-TextValue $ws 'D22' '186.75'
Set-TextValue $ws 'E22' '  -11.13%  '
Set-TextValue $ws 'D23' '6.256'
Set-TextValue $ws 'E23' '  -8.12%  '
Set-TextValue $ws 'D24' '1.006'
Set-TextValue $ws 'E24' '  +0.20%  '
Set-TextValue $ws 'D25' '144.88'
Set-TextValue $ws 'E25' '  -5.21%  '
Set-TextValue $ws 'D26' '7.488'
Set-TextValue $ws 'E26' '  -7.38%  '
Set-TextValue $ws 'D27' '0.1174'
Set-TextValue $ws 'E27' '  -6.74%  '
Set-TextValue $ws 'D28' '15.80'
Set-TextValue $ws 'E28' '  -3.67%  '
Set-TextValue $ws 'D29' '1.342'
Set-TextValue $ws 'E29' '  -5.17%  '
Set-TextValue $ws 'E30' '  -6.51%  '
Set-TextValue $ws 'D31' '1.342'
Set-TextValue $ws 'E31' '  -5.61%  '
Set-TextValue $ws 'E32' '  -6.92%  '
Set-TextValue $ws 'D33' '3.508'
Set-TextValue $ws 'E33' '  -8.01%  '
Set-TextValue $ws 'E34' '  -5.40%  '
Set-TextValue $ws 'E35' '  -2.60%  '
Set-TextValue $ws 'D36' '0.5995'
Set-TextValue $ws 'E36' '  -6.44%  '
Set-TextValue $ws 'D37' '2.364'
Set-TextValue $ws 'E37' '  -5.49%  '
Set-TextValue $ws 'D38' '2.682'
Set-TextValue $ws 'E38' '  -1.21%  '
Set-TextValue $ws 'E39' '  -4.43%  '
Set-TextValue $ws 'D40' '1.090.99'
Set-TextValue $ws 'E40' '  -4.49%  '
Set-TextValue $ws 'D41' '0.8617'
Set-TextValue $ws 'E41' '  -1.99%  '
Set-TextValue $ws 'D42' '5.830'
Set-TextValue $ws 'E42' '  -9.27%  '
Set-TextValue $ws 'E43' '  -0.05%  '
Set-TextValue $ws 'D44' '99.54'
Set-TextValue $ws 'E44' '  -0.59%  '
Set-TextValue $ws 'D45' '1.850.55'
Set-TextValue $ws 'E45' '  -4.94%  '
Set-TextValue $ws 'E46' '  +6.24%  '
Set-TextValue $ws 'D47' '56.45'
Set-TextValue $ws 'E47' '  -5.72%  '
Set-TextValue $ws 'E48' '  +0.45%  '
Set-TextValue $ws 'D49' '8.096'
Set-TextValue $ws 'E49' '  -3.25%  '
Set-TextValue $ws 'B50' 'Cronos'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws 'D50' '0.05234'
Set-TextValue $ws 'E50' '  -4.61%  '
Set-TextValue $ws 'B51' 'Mantle'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws 'D51' '0.4323'
Set-TextValue $ws 'E51' '  -3.41%  '
